$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13) for every data row
# (rows 2 through 525) as part of an automatic daily data refresh.
$ws.Range("C2:C525").Value2 = 45182
